# Swap the "First Name" / "Last Name" columns (A <-> B) on the names sheet
# so Last Name precedes First Name, per the commit "change last name before first name".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Swap the column widths (A <-> B) -----------------------------------
# Read the exact current widths via the pixel-width getter (full precision),
# convert pixels -> stored character width, then feed the cross value back
# into ColumnWidth (which itself re-derives the stored width).
$pxA = $ws.Columns.Item(1).Width
$pxB = $ws.Columns.Item(2).Width
$storedWidthA = ($pxA - 3.75) / 5.25
$storedWidthB = ($pxB - 3.75) / 5.25

$ws.Columns.Item(1).ColumnWidth = $storedWidthB - (5.0 / 6.0)
$ws.Columns.Item(2).ColumnWidth = $storedWidthA - (5.0 / 6.0)

# --- 2. Swap the cell contents of columns A and B, row by row --------------
$lastRow = 291
for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $textA = $cellA.Text
    $textB = $cellB.Text
    $cellA.Value = $textB
    $cellB.Value = $textA
}

# --- 3. Update the AutoFilter range from B1:C291 to A1:B291 -----------------
$ws.AutoFilterMode = $false
$ws.Range("A1:B291").AutoFilter()

# --- 4. Update the workbook-level _FilterDatabase defined name -------------
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -like "*_FilterDatabase*") {
        $nm.RefersTo = "='" + $ws.Name + "'!`$A`$1:`$B`$291"
    }
}

# --- 5. Select column A (matches the resulting selection in the saved file)-
$ws.Range("A:A").Select()
